$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24, shifting existing rows 24-39 down to 25-40
$ws.Rows.Item(24).Insert()

# Populate the new row 24 with this week's data (a new Black Amber / Primera
# record for Región Metropolitana), matching the other rows' layout.
$row = 24
$ws.Cells.Item($row, 1).Value = 11
$ws.Cells.Item($row, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item($row, 3).Value = "Bíobío"
$ws.Cells.Item($row, 4).Value = 44566
$ws.Cells.Item($row, 5).Value = 8
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100103
$ws.Cells.Item($row, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item($row, 9).Value = 100103002
$ws.Cells.Item($row, 10).Value = "Ciruela"
$ws.Cells.Item($row, 11).Value = "Black Amber"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 200
$ws.Cells.Item($row, 14).Value = 14000
$ws.Cells.Item($row, 15).Value = 15000
$ws.Cells.Item($row, 16).Value = 14400
$ws.Cells.Item($row, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item($row, 18).Value = "Región Metropolitana"
$ws.Cells.Item($row, 19).Value = 800
$ws.Cells.Item($row, 20).Value = 18

# Match the date cell format used by the other "Fecha" cells in column D
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
